$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace backslash-style relative paths with forward-slash style paths
$ws.Range("B2").Value = "./Resources"
$ws.Range("B3").Value = "./Resources"
$ws.Range("B4").Value = "./Resources"
$ws.Range("B5").Value = "./Resources"
$ws.Range("B6").Value = "./Resources/G_Money"
$ws.Range("B7").Value = "./Resources/yeshut"
$ws.Range("B8").Value = "./Resources/Changemat"
$ws.Range("B9").Value = "./Resources/cox"
$ws.Range("B10").Value = "./Resources/yeshut_xc"

# Update the view selection (also resets topLeftCell to the default)
$ws.Range("I9").Select()
